$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "ea42df71-985a-4d82-8ca7-ae60a86de51d.md"
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("D2").Value = "2016-36-19 02:36:44"
$wsOverview.Range("A3").Value = "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-37-19 02:37:36"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "ea42df71-985a-4d82-8ca7-ae60a86de51d.md"
$wsZhCn.Range("D2").Value = "ea42df71-985a-4d82-8ca7-ae60a86de51d.d82585e3266a853a34f383049954b70dfab256ed.zh-cn.xlf"
$wsZhCn.Range("F2").Value = "ea42df71-985a-4d82-8ca7-ae60a86de51d.md"
$wsZhCn.Range("G2").Value = "ea42df71-985a-4d82-8ca7-ae60a86de51d.d82585e3266a853a34f383049954b70dfab256ed.zh-cn.xlf"
$wsZhCn.Range("A3").Value = "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.6a0d9d2bb57d903e710272a8ba7ce0f4e005fb66.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-19 02:37:32"
$wsZhCn.Range("F3").Value = "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md"
$wsZhCn.Range("G3").Value = "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.6a0d9d2bb57d903e710272a8ba7ce0f4e005fb66.zh-cn.xlf"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "ea42df71-985a-4d82-8ca7-ae60a86de51d.md"
$wsDeDe.Range("D2").Value = "ea42df71-985a-4d82-8ca7-ae60a86de51d.d82585e3266a853a34f383049954b70dfab256ed.de-de.xlf"
$wsDeDe.Range("F2").Value = "ea42df71-985a-4d82-8ca7-ae60a86de51d.md"
$wsDeDe.Range("G2").Value = "ea42df71-985a-4d82-8ca7-ae60a86de51d.d82585e3266a853a34f383049954b70dfab256ed.de-de.xlf"
$wsDeDe.Range("A3").Value = "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.6a0d9d2bb57d903e710272a8ba7ce0f4e005fb66.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-19 02:37:36"
$wsDeDe.Range("F3").Value = "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md"
$wsDeDe.Range("G3").Value = "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.6a0d9d2bb57d903e710272a8ba7ce0f4e005fb66.de-de.xlf"
